# Add six new status-report rows (A80:C85) to Sheet1, matching the
# "hours update src update" commit.
#
# Column A uses the same date-formatted style (numFmtId 14) as the rest of
# the date column, so copy the format from an existing date cell before
# filling in new values rather than setting NumberFormat directly (which
# would create a brand-new custom style/number-format entry instead of
# reusing the existing one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

[void]$ws.Range("A79").Copy()
[void]$ws.Range("A80:A85").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A80").Value = 40275
$ws.Range("B80").Value = 6
$ws.Range("C80").Value = "beagleboard avr-can interface"

$ws.Range("A81").Value = 40277
$ws.Range("B81").Value = 1
$ws.Range("C81").Value = "buying parts for interface board"

$ws.Range("A82").Value = 40277
$ws.Range("B82").Value = 6.5
$ws.Range("C82").Value = "interface board soldering and assembly"

$ws.Range("A83").Value = 40278
$ws.Range("B83").Value = 8
$ws.Range("C83").Value = "interface board debug"

$ws.Range("A84").Value = 40279
$ws.Range("B84").Value = 3
$ws.Range("C84").Value = "pair programming with Erica"

$ws.Range("A85").Value = 40279
$ws.Range("B85").Value = 6
$ws.Range("C85").Value = "serial interface debug"

# Move the viewport / selection to match the author's final cursor position.
[void]$ws.Range("C86").Select()
